$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: RED won with killer word selected
$ws.Range("C24").Value = "RED"
$ws.Range("D24").Value = 2.333333333333333
$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 0.3333333333333333
$ws.Range("G24").Value = 0.3333333333333333
$ws.Range("H24").Value = "killer word selected"
$ws.Range("I24").Value = 3
$ws.Range("J24").Value = 3
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0

# Row 25: BLUE won with killer word selected
$ws.Range("C25").Value = "BLUE"
$ws.Range("D25").Value = 2.25
$ws.Range("E25").Value = 2.333333333333333
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 0.6666666666666666
$ws.Range("H25").Value = "killer word selected"
$ws.Range("I25").Value = 4
$ws.Range("J25").Value = 3
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 1
